$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.041.72"
$ws.Range("E2").Value = "  +2.31%  "
$ws.Range("D3").Value = "'2.458.83"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'577.04"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").Value = "'146.69"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("D9").Value = "'2.458.06"
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("E11").Value = "  +1.81%  "
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").Value = "'0.354"
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("E14").Value = "  +8.80%  "
$ws.Range("E15").Value = "  +3.49%  "
$ws.Range("D16").Value = "'2.904.34"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D17").Value = "'62.940.83"
$ws.Range("E17").Value = "  +2.84%  "
$ws.Range("D18").Value = "'2.464.82"
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("D19").Value = "'7.98"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "'11.12"
$ws.Range("E20").Value = "  +4.30%  "
$ws.Range("D21").Value = "'330.04"
$ws.Range("E21").Value = "  +1.76%  "
$ws.Range("D22").Value = "'2.23"
$ws.Range("E22").Value = "  +12.93%  "
$ws.Range("D23").Value = "'4.13"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D25").Value = "'66.47"
$ws.Range("E25").Value = "  +2.14%  "
$ws.Range("D26").Value = "'662.13"
$ws.Range("E26").Value = "  +7.39%  "
$ws.Range("E27").Value = "  +15.63%  "
$ws.Range("D28").Value = "'8.83"
$ws.Range("E28").Value = "  +6.21%  "
$ws.Range("E29").Value = "  +4.57%  "
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("D32").Value = "'1.44"
$ws.Range("E32").Value = "  +3.37%  "
$ws.Range("E33").Value = "  +5.14%  "
$ws.Range("D34").Value = "'0.139"
$ws.Range("E34").Value = "  +4.30%  "
$ws.Range("E35").Value = "  +4.30%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  +3.32%  "
$ws.Range("E38").Value = "  +3.58%  "
$ws.Range("D39").Value = "'153.08"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  +2.54%  "
$ws.Range("D42").Value = "'0.0₆0352"
$ws.Range("E42").Value = "  +23.55%  "
$ws.Range("E43").Value = "  +5.87%  "
$ws.Range("E44").Value = "  +3.37%  "
$ws.Range("D45").Value = "'42.33"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "'15.14"
$ws.Range("E47").Value = "  +28.10%  "
$ws.Range("D48").Value = "'146.52"
$ws.Range("E48").Value = "  +2.42%  "
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("E50").Value = "  +3.72%  "
$ws.Range("D51").Value = "'0.608"
$ws.Range("E51").Value = "  +2.07%  "
